$wb = $excel.ActiveWorkbook

# Map of per-sheet cell updates: column letter -> value (market-data refresh)
# Cells listed under "Clear" are removed entirely (no cached price available),
# matching the source diff which drops those <c> elements.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3685.5715
$ws.Range("J40").Value = 3869.4348
$ws.Range("L40").Value = 3869.4348
$ws.Range("N40").Value = -4219.4348

$ws.Range("H53").Value = 1081
$ws.Range("I53").Value = 430.66666
$ws.Range("J53").Value = 1861.4
$ws.Range("K53").Value = 430.66666
$ws.Range("L53").Value = 1861.4
$ws.Range("M53").Value = 206.33334
$ws.Range("N53").Value = -3135.4

$ws.Range("H74").Value = 7544.0713
$ws.Range("I74").Value = 7355.154
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 7355.154
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -6419.154
$ws.Range("N74").Value = -11872

$ws.Range("H77").Value = 7544.0713
$ws.Range("I77").Value = 7355.154
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 36775.77
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -32095.77
$ws.Range("N77").Value = -59360

$ws.Range("H100").Value = 3057.7144
$ws.Range("I100").Value = 1978.6666
$ws.Range("K100").Value = 1978.6666
$ws.Range("M100").Value = -1437.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2130.6843
$ws.Range("I102").Value = 1940.1765
$ws.Range("K102").Value = 1940.1765
$ws.Range("M102").Value = -318.1765

$ws.Range("H110").Value = 1518.8
$ws.Range("I110").Value = 1256.8334
$ws.Range("K110").Value = 1256.8334
$ws.Range("M110").Value = 788.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 39278.93
$ws.Range("I20").Value = 58378.445
$ws.Range("K20").Value = 58378.445
$ws.Range("M20").Value = -58131.445

$ws.Range("H64").Value = 1022.25
$ws.Range("J64").Value = 1022.25
$ws.Range("L64").Value = 1022.25
$ws.Range("N64").Value = -1472.25

$ws.Range("H67").Value = 1022.25
$ws.Range("J67").Value = 1022.25
$ws.Range("L67").Value = 1022.25
$ws.Range("N67").Value = -2582.25

$ws.Range("H128").Value = 7008.8184
$ws.Range("I128").Value = 7008.8184
$ws.Range("K128").Value = 21026.4552
$ws.Range("M128").Value = -18536.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 173995
$ws.Range("J118").Value = 173995
$ws.Range("L118").Value = 173995
$ws.Range("N118").Value = -177309

$ws.Range("H122").Value = 3455.628
$ws.Range("I122").Value = 2985.2964
$ws.Range("J122").Value = 4249.3125
$ws.Range("K122").Value = 8955.889200000001
$ws.Range("L122").Value = 12747.9375
$ws.Range("M122").Value = -6505.889200000001
$ws.Range("N122").Value = -17647.9375

$ws.Range("H132").Value = 4645.143
$ws.Range("I132").Value = 5219.7144
$ws.Range("J132").Value = 4070.5715
$ws.Range("K132").Value = 15659.1432
$ws.Range("L132").Value = 12211.7145
$ws.Range("M132").Value = -13129.1432
$ws.Range("N132").Value = -17271.7145

$ws.Range("H134").Value = 2421.8572
$ws.Range("I134").Value = 2421.8572
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7265.571599999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4730.571599999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3599.077
$ws.Range("J39").Value = 3599.077
$ws.Range("L39").Value = 10797.231
$ws.Range("N39").Value = -11385.231

$ws.Range("H55").Value = 2705.4167
$ws.Range("I55").Value = 1225
$ws.Range("J55").Value = 3001.5
$ws.Range("K55").Value = 3675
$ws.Range("L55").Value = 9004.5
$ws.Range("M55").Value = -3498
$ws.Range("N55").Value = -9358.5

$ws.Range("H80").Value = 2110.5
$ws.Range("J80").Value = 2110.5
$ws.Range("L80").Value = 6331.5
$ws.Range("N80").Value = -8203.5

$ws.Range("H81").Value = 2833
$ws.Range("I81").Value = 2475.7144
$ws.Range("K81").Value = 7427.1432
$ws.Range("M81").Value = -6304.1432

$ws.Range("H83").Value = 2110.5
$ws.Range("J83").Value = 2110.5
$ws.Range("L83").Value = 18994.5
$ws.Range("N83").Value = -28354.5

$ws.Range("H84").Value = 2833
$ws.Range("I84").Value = 2475.7144
$ws.Range("K84").Value = 22281.4296
$ws.Range("M84").Value = -16665.4296

$ws.Range("H113").Value = 1615.6666
$ws.Range("J113").Value = 1778.8
$ws.Range("L113").Value = 5336.4
$ws.Range("N113").Value = -9676.4

$ws.Range("H118").Value = 3000
$ws.Range("I118").Value = 2000
$ws.Range("J118").Value = 4000
$ws.Range("K118").Value = 6000
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = -4757
$ws.Range("N118").Value = -14486

$ws.Range("H131").Value = 1704.1282
$ws.Range("J131").Value = 1765.3636
$ws.Range("L131").Value = 5296.0908
$ws.Range("N131").Value = -15376.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 29888
$ws.Range("J93").Value = 29888
$ws.Range("L93").Value = 29888
$ws.Range("N93").Value = -33632

$ws.Range("H102").Value = 2825.2727
$ws.Range("I102").Value = 2885.5
$ws.Range("J102").Value = 2664.6667
$ws.Range("K102").Value = 2885.5
$ws.Range("L102").Value = 2664.6667
$ws.Range("M102").Value = -1263.5
$ws.Range("N102").Value = -5908.6667

$ws.Range("H113").Value = 14611
$ws.Range("J113").Value = 11420
$ws.Range("L113").Value = 11420
$ws.Range("N113").Value = -15760

$ws.Range("H122").Value = 2151.3333
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2151.3333
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6453.999899999999
$ws.Range("N122").Value = -11353.9999
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 3783.1667
$ws.Range("I126").Value = 3565.8333
$ws.Range("J126").Value = 4000.5
$ws.Range("K126").Value = 10697.4999
$ws.Range("L126").Value = 12001.5
$ws.Range("M126").Value = -8227.499899999999
$ws.Range("N126").Value = -16941.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5708.2607
$ws.Range("I7").Value = 5370.5884
$ws.Range("J7").Value = 6665
$ws.Range("K7").Value = 5370.5884
$ws.Range("L7").Value = 6665
$ws.Range("M7").Value = -5258.5884
$ws.Range("N7").Value = -6889

$ws.Range("H40").Value = 5300.5557
$ws.Range("I40").Value = 3992.5
$ws.Range("J40").Value = 7916.6665
$ws.Range("K40").Value = 3992.5
$ws.Range("L40").Value = 7916.6665
$ws.Range("M40").Value = -3856.5
$ws.Range("N40").Value = -8188.6665

$ws.Range("H46").Value = 13951.6
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 13951.6
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 13951.6
$ws.Range("N46").Value = -14327.6
$ws.Range("M46").ClearContents()

$ws.Range("H61").Value = 1424.8182
$ws.Range("I61").Value = 1329.5555
$ws.Range("K61").Value = 1329.5555
$ws.Range("M61").Value = -1127.5555

$ws.Range("H68").Value = 3887.125
$ws.Range("I68").Value = 5019.6
$ws.Range("K68").Value = 5019.6
$ws.Range("M68").Value = -4270.6

$ws.Range("H71").Value = 3887.125
$ws.Range("I71").Value = 5019.6
$ws.Range("K71").Value = 25098
$ws.Range("M71").Value = -21354

$ws.Range("H113").Value = 1424.8182
$ws.Range("I113").Value = 1329.5555
$ws.Range("K113").Value = 1329.5555
$ws.Range("M113").Value = 840.4445000000001

$ws.Range("H126").Value = 5708.2607
$ws.Range("I126").Value = 5370.5884
$ws.Range("J126").Value = 6665
$ws.Range("K126").Value = 16111.7652
$ws.Range("L126").Value = 19995
$ws.Range("M126").Value = -13641.7652
$ws.Range("N126").Value = -24935

$ws.Range("H132").Value = 4907.92
$ws.Range("I132").Value = 4664.7144
$ws.Range("J132").Value = 6184.75
$ws.Range("K132").Value = 13994.1432
$ws.Range("L132").Value = 18554.25
$ws.Range("M132").Value = -11464.1432
$ws.Range("N132").Value = -23614.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 155000
$ws.Range("J57").Value = 155000
$ws.Range("L57").Value = 155000
$ws.Range("N57").Value = -156508

$ws.Range("H81").Value = 3622.611
$ws.Range("I81").Value = 2714.1
$ws.Range("J81").Value = 4758.25
$ws.Range("K81").Value = 5428.2
$ws.Range("L81").Value = 9516.5
$ws.Range("M81").Value = -4367.2
$ws.Range("N81").Value = -11638.5

$ws.Range("H84").Value = 3622.611
$ws.Range("I84").Value = 2714.1
$ws.Range("J84").Value = 4758.25
$ws.Range("K84").Value = 27141
$ws.Range("L84").Value = 47582.5
$ws.Range("M84").Value = -21837
$ws.Range("N84").Value = -58190.5

$ws.Range("H103").Value = 112995
$ws.Range("J103").Value = 112995
$ws.Range("L103").Value = 112995
$ws.Range("N103").Value = -115339

$ws.Range("H132").Value = 2212.0967
$ws.Range("I132").Value = 2191.7307
$ws.Range("K132").Value = 6575.1921
$ws.Range("M132").Value = -4045.1921

Write-Output "Updated 43 leve-profit rows across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
